$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1128.1818
$ws.Range("I17").Value = 684.2857
$ws.Range("J17").Value = 1335.3334
$ws.Range("K17").Value = 2052.8571
$ws.Range("L17").Value = 4006.0002
$ws.Range("M17").Value = -1884.8571
$ws.Range("N17").Value = -4342.0002
$ws.Range("H28").Value = 4078.2
$ws.Range("I28").Value = 311.66666
$ws.Range("J28").Value = 9728
$ws.Range("K28").Value = 311.66666
$ws.Range("L28").Value = 9728
$ws.Range("M28").Value = 173.33334
$ws.Range("N28").Value = -10698
$ws.Range("H62").Value = 5426.25
$ws.Range("I62").Value = 5426.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5426.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4802.25
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 45000
$ws.Range("L63").Value = 45000
$ws.Range("N63").Value = -46248
$ws.Range("H65").Value = 5426.25
$ws.Range("I65").Value = 5426.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 27131.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -24011.25
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 45000
$ws.Range("J66").Value = 45000
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141240
$ws.Range("J75").Value = 21662.8
$ws.Range("L75").Value = 21662.8
$ws.Range("N75").Value = -23534.8
$ws.Range("H76").Value = 65187.688
$ws.Range("I76").Value = 65187.688
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 65187.688
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -64872.688
$ws.Range("N76").ClearContents()
$ws.Range("J78").Value = 21662.8
$ws.Range("L78").Value = 64988.39999999999
$ws.Range("N78").Value = -74348.39999999999
$ws.Range("H79").Value = 65187.688
$ws.Range("I79").Value = 65187.688
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 65187.688
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -64095.688
$ws.Range("N79").ClearContents()
$ws.Range("H106").Value = 62626172
$ws.Range("I106").Value = 144198.42
$ws.Range("J106").Value = 500000000
$ws.Range("K106").Value = 144198.42
$ws.Range("L106").Value = 500000000
$ws.Range("M106").Value = -143567.42
$ws.Range("N106").Value = -500001262
$ws.Range("H132").Value = 8404780
$ws.Range("I132").Value = 10205625
$ws.Range("J132").Value = 835
$ws.Range("K132").Value = 30616875
$ws.Range("L132").Value = 2505
$ws.Range("M132").Value = -30614345
$ws.Range("N132").Value = -7565
$ws.Range("H137").Value = 976.24243
$ws.Range("I137").Value = 983.21875
$ws.Range("J137").Value = 753
$ws.Range("K137").Value = 2949.65625
$ws.Range("L137").Value = 2259
$ws.Range("M137").Value = -399.65625
$ws.Range("N137").Value = -7359

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14446.583
$ws.Range("I32").Value = 15963.014
$ws.Range("K32").Value = 15963.014
$ws.Range("M32").Value = -15676.014
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 41933
$ws.Range("J64").Value = 41933
$ws.Range("L64").Value = 41933
$ws.Range("N64").Value = -42429
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 41933
$ws.Range("J67").Value = 41933
$ws.Range("L67").Value = 41933
$ws.Range("N67").Value = -43649
$ws.Range("H96").Value = 24614.666
$ws.Range("J96").Value = 24614.666
$ws.Range("L96").Value = 24614.666
$ws.Range("N96").Value = -30106.666
$ws.Range("H97").Value = 408
$ws.Range("I97").Value = 326.66666
$ws.Range("K97").Value = 326.66666
$ws.Range("M97").Value = 169.33334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2401.4849
$ws.Range("I20").Value = 2568.6667
$ws.Range("J20").Value = 2200.8667
$ws.Range("K20").Value = 2568.6667
$ws.Range("L20").Value = 2200.8667
$ws.Range("M20").Value = -2321.6667
$ws.Range("N20").Value = -2694.8667
$ws.Range("H99").Value = 761.9524
$ws.Range("I99").Value = 584.2857
$ws.Range("J99").Value = 850.7857
$ws.Range("K99").Value = 584.2857
$ws.Range("L99").Value = 850.7857
$ws.Range("M99").Value = 913.7143
$ws.Range("N99").Value = -3846.7857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1120.375
$ws.Range("I16").Value = 1025
$ws.Range("K16").Value = 1025
$ws.Range("M16").Value = -738
$ws.Range("H31").Value = 11114851
$ws.Range("I31").Value = 3480
$ws.Range("J31").Value = 23813562
$ws.Range("K31").Value = 3480
$ws.Range("L31").Value = 23813562
$ws.Range("M31").Value = -3185
$ws.Range("N31").Value = -23814152
$ws.Range("H34").Value = 11114851
$ws.Range("I34").Value = 3480
$ws.Range("J34").Value = 23813562
$ws.Range("K34").Value = 3480
$ws.Range("L34").Value = 23813562
$ws.Range("M34").Value = -3278
$ws.Range("N34").Value = -23813966
$ws.Range("H113").Value = 1120.375
$ws.Range("I113").Value = 1025
$ws.Range("K113").Value = 1025
$ws.Range("M113").Value = 1145
$ws.Range("H122").Value = 599.63635
$ws.Range("I122").Value = 610.375
$ws.Range("J122").Value = 571
$ws.Range("K122").Value = 1831.125
$ws.Range("L122").Value = 1713
$ws.Range("M122").Value = 618.875
$ws.Range("N122").Value = -6613
$ws.Range("H132").Value = 2138.1667
$ws.Range("I132").Value = 1245.2693
$ws.Range("J132").Value = 4459.7
$ws.Range("K132").Value = 3735.8079
$ws.Range("L132").Value = 13379.1
$ws.Range("M132").Value = -1205.8079
$ws.Range("N132").Value = -18439.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 141.09525
$ws.Range("I12").Value = 152.5
$ws.Range("J12").Value = 136.53334
$ws.Range("K12").Value = 457.5
$ws.Range("L12").Value = 409.60002
$ws.Range("M12").Value = -284.5
$ws.Range("N12").Value = -755.6000200000001
$ws.Range("H131").Value = 7683.933
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 50757.5
$ws.Range("J34").Value = 50757.5
$ws.Range("L34").Value = 50757.5
$ws.Range("N34").Value = -51293.5
$ws.Range("H70").Value = 17003920
$ws.Range("I70").Value = 18892778
$ws.Range("K70").Value = 18892778
$ws.Range("M70").Value = -18892508
$ws.Range("H73").Value = 17003920
$ws.Range("I73").Value = 18892778
$ws.Range("K73").Value = 18892778
$ws.Range("M73").Value = -18891842
$ws.Range("H76").Value = 50757.5
$ws.Range("J76").Value = 50757.5
$ws.Range("L76").Value = 50757.5
$ws.Range("N76").Value = -51387.5
$ws.Range("H79").Value = 50757.5
$ws.Range("J79").Value = 50757.5
$ws.Range("L79").Value = 50757.5
$ws.Range("N79").Value = -52941.5
$ws.Range("H97").Value = 974.3889
$ws.Range("I97").Value = 977.6923
$ws.Range("J97").Value = 965.8
$ws.Range("K97").Value = 977.6923
$ws.Range("L97").Value = 965.8
$ws.Range("M97").Value = -481.6923
$ws.Range("N97").Value = -1957.8
$ws.Range("H132").Value = 24472.512
$ws.Range("I132").Value = 31259.912
$ws.Range("J132").Value = 3493.2727
$ws.Range("K132").Value = 93779.736
$ws.Range("L132").Value = 10479.8181
$ws.Range("M132").Value = -91249.736
$ws.Range("N132").Value = -15539.8181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15450
$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16560
$ws.Range("H76").Value = 17144
$ws.Range("J76").Value = 17144
$ws.Range("L76").Value = 17144
$ws.Range("N76").Value = -17820
$ws.Range("H79").Value = 17144
$ws.Range("J79").Value = 17144
$ws.Range("L79").Value = 17144
$ws.Range("N79").Value = -19484
$ws.Range("H122").Value = 2950.3684
$ws.Range("I122").Value = 3673
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 11019
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -8569
$ws.Range("N122").Value = -11800
$ws.Range("H132").Value = 1978.8214
$ws.Range("I132").Value = 2404
$ws.Range("J132").Value = 1703.7059
$ws.Range("K132").Value = 7212
$ws.Range("L132").Value = 5111.1177
$ws.Range("M132").Value = -4682
$ws.Range("N132").Value = -10171.1177
$ws.Range("H133").Value = 19225.111
$ws.Range("J133").Value = 19225.111
$ws.Range("L133").Value = 19225.111
$ws.Range("N133").Value = -24285.111

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2872.5908
$ws.Range("I132").Value = 2224.0833
$ws.Range("J132").Value = 3650.8
$ws.Range("K132").Value = 6672.249899999999
$ws.Range("L132").Value = 10952.4
$ws.Range("M132").Value = -4142.249899999999
$ws.Range("N132").Value = -16012.4
